# Update layer definitions to use names from Staurenghi et al. 2014.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (surface_id 6): ONL-MZ -> ONL, description updated
$ws.Range("B7").Value = "ONL"
$ws.Range("C7").Value = "Henle's fiber layer, outer nuclear layer, & myoid zone"

# Row 8 (surface_id 7): EZ description updated
$ws.Range("B8").Value = "EZ"
$ws.Range("C8").Value = "ellipsoid zone"

# Row 9 (surface_id 8): IZOS -> OS
$ws.Range("B9").Value = "OS"
$ws.Range("C9").Value = "outer segments"

# Row 10 (surface_id 9): OZOS -> IZ
$ws.Range("B10").Value = "IZ"
$ws.Range("C10").Value = "interdigitation zone"

# Row 11: surface_id 10 -> 11, RPE description updated
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "RPE"
$ws.Range("C11").Value = "RPE/Bruch's complex"

# Row 12: surface_id 11 -> 12, outer RPE -> CHOR-SCL
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "CHOR-SCL"
$ws.Range("C12").Value = "choroid-sclera"

# Column C is now wider to fit the longer descriptions
$ws.Columns("C").ColumnWidth = 44.7

# Selection moved to E11
$ws.Range("E11").Select()
